# Add new translation keys/values (filter_*, expand_all) to the
# "translations" sheet, rows 64-68, columns A (key) and B (German value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

$rows = @(
    @("filter_involved", "Person"),
    @("filter_creator", "Autor"),
    @("filter_journal", "Zeitschrift"),
    @("filter_type", "Textart"),
    @("expand_all", "alle ausklappen")
)

$startRow = 64
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $pair = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

$ws.Range("A69").Select()
